$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Desired Freq" value (B26); dependent formulas B27/B29 recalc automatically.
$ws.Range("B26").Value = 1000

# Move/update the active selection to B24.
$ws.Range("B24").Select()
